# "new selector prev month"
# Shift the two date cells back to the previous month/year pair and move the
# active selection from B3 up to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# fechaInicial (A2) / FechaFinal (A3) date values — pushed back to an earlier period
$ws.Range("B2").Value = 43831
$ws.Range("B3").Value = 43846

# Move the active cell/selection from B3 to B2
$ws.Range("B2").Select()
